$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to match the scraped player's name
$ws.Name = "Deepak Chahar"

# Insert a new "matchNo" column before column A, shifting teamName..result right
$ws.Columns("A").Insert()

# Insert a new row before the existing data row, shifting it down (row2 -> row3)
$ws.Rows(2).Insert()

# Force text storage for the numeric-looking stat columns so "0"/"1"/"100.00"
# are kept as text rather than being auto-converted to numbers (matches the
# source data, which stores every field - including numbers - as text).
$ws.Range("E2:I2").NumberFormat = "@"
$ws.Range("E3:H3").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# New row 2: 38th match vs Kolkata Knight Riders
$ws.Range("A2").Value = "38th"
$ws.Range("B2").Value = "Chennai Super Kings"
$ws.Range("C2").Value = "Deepak Chahar"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "100.00"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Abu Dhabi"
$ws.Range("L2").Value = "September 26"
$ws.Range("M2").Value = "Super Kings won by 2 wickets"

# Row 3 (previously row 2): 12th match vs Rajasthan Royals, now with matchNo filled in
$ws.Range("A3").Value = "12th"
$ws.Range("B3").Value = "Chennai Super Kings"
$ws.Range("C3").Value = "Deepak Chahar"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "0"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "-"
$ws.Range("J3").Value = "Rajasthan Royals"
$ws.Range("K3").Value = "Wankhede"
$ws.Range("L3").Value = "April 19"
$ws.Range("M3").Value = "Super Kings won by 45 runs"
